$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 274
$ws.Range("F5").Value = 2028
$ws.Range("F6").Value = 82
$ws.Range("F7").Value = 513
$ws.Range("F8").Value = 432
$ws.Range("F9").Value = 220
$ws.Range("F10").Value = 7398
$ws.Range("F11").Value = 207
$ws.Range("F12").Value = 561
$ws.Range("F13").Value = 773
$ws.Range("F14").Value = 73
$ws.Range("F15").Value = 3159
$ws.Range("F16").Value = 1818
$ws.Range("F17").Value = 161
$ws.Range("F18").Value = 13
$ws.Range("F19").Value = 53
$ws.Range("F20").Value = 112
$ws.Range("F21").Value = 177
$ws.Range("F22").Value = 125
$ws.Range("F24").Value = 185
$ws.Range("F25").Value = 83
$ws.Range("F26").Value = 988
$ws.Range("F27").Value = 212
$ws.Range("F28").Value = 4141

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 80
$ws.Range("F3").Value = 30
$ws.Range("F4").Value = 15

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 131
$ws.Range("F3").Value = 735

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 131
$ws.Range("F4").Value = 735
$ws.Range("F6").Value = 80
$ws.Range("F7").Value = 274
$ws.Range("F8").Value = 2028
$ws.Range("F9").Value = 30
$ws.Range("F10").Value = 15
$ws.Range("F11").Value = 82
$ws.Range("F12").Value = 513
$ws.Range("F13").Value = 432
$ws.Range("F14").Value = 220
$ws.Range("F15").Value = 7398
$ws.Range("F16").Value = 207
$ws.Range("F17").Value = 561
$ws.Range("F18").Value = 773
$ws.Range("F19").Value = 73
$ws.Range("F20").Value = 3159
$ws.Range("F21").Value = 1818
$ws.Range("F22").Value = 161
$ws.Range("F23").Value = 13
$ws.Range("F24").Value = 53
$ws.Range("F25").Value = 112
$ws.Range("F26").Value = 177
$ws.Range("F27").Value = 125
$ws.Range("F29").Value = 185
$ws.Range("F30").Value = 83
$ws.Range("F31").Value = 988
$ws.Range("F32").Value = 212
$ws.Range("F33").Value = 4141
